# Weekly driver report update for 2025-05-05
# Updates the "Driver Summary" sheet of the ABM driver summary workbook with
# refreshed roaming-impact figures (client counts, critical/warning/good
# minute sums, totals) and re-sorts the "Good Drivers" driver-version rows
# (14-21) to reflect the new vintage ordering produced by this week's run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" summary block (rows 3-4) ---
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 218
$ws.Range("D3").Value = 98.59999999999999

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 218

# --- "Good Drivers" detail block (rows 12-24) ---

# Row 13: Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.0.4 (unchanged position)
$ws.Range("B13").Value = 1293197
$ws.Range("C13").Value = 4322
$ws.Range("E13").Value = 1990
$ws.Range("F13").Value = 1298528

# Row 14: now Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4 (moved up from row 16)
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4"
$ws.Range("B14").Value = 31517
$ws.Range("C14").Value = 112
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 66
$ws.Range("F14").Value = 31629
$ws.Range("H14").Value = "22.220.0.4"
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = "2023-03-28"
$ws.Range("J14").ClearFormats()

# Row 15: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8 (moved up from row 14)
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8"
$ws.Range("B15").Value = 467311
$ws.Range("C15").Value = 772
$ws.Range("D15").Value = 931
$ws.Range("E15").Value = 706
$ws.Range("F15").Value = 469014
$ws.Range("H15").Value = "23.10.0.8"
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = "2023-10-30"
$ws.Range("J15").ClearFormats()

# Row 16: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3 (moved down from row 15)
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3"
$ws.Range("B16").Value = 455081
$ws.Range("C16").Value = 1861
$ws.Range("D16").Value = 52
$ws.Range("E16").Value = 639
$ws.Range("F16").Value = 456994
$ws.Range("H16").Value = "23.120.0.3"
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "2025-02-05"
$ws.Range("J16").ClearFormats()

# Row 17: Intel(R) Wi-Fi 6E AX211 160MHz - 22.230.0.8 (unchanged position)
$ws.Range("B17").Value = 1787924
$ws.Range("C17").Value = 3326
$ws.Range("E17").Value = 3038
$ws.Range("F17").Value = 1793864

# Row 18: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3 (moved up from row 19)
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3"
$ws.Range("B18").Value = 218767
$ws.Range("C18").Value = 334
$ws.Range("D18").Value = 313
$ws.Range("E18").Value = 573
$ws.Range("F18").Value = 219414
$ws.Range("H18").Value = "23.70.2.3"
$ws.Range("J18").NumberFormat = "@"
$ws.Range("J18").Value = "2024-07-23"
$ws.Range("J18").ClearFormats()

# Row 19: now Intel(R) Wi-Fi 6E AX211 160MHz - 22.110.1.1 (moved down from row 18)
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.110.1.1"
$ws.Range("B19").Value = 135467
$ws.Range("C19").Value = 189
$ws.Range("D19").Value = 263
$ws.Range("E19").Value = 196
$ws.Range("F19").Value = 135919
$ws.Range("H19").Value = "22.110.1.1"
$ws.Range("J19").NumberFormat = "@"
$ws.Range("J19").Value = "2022-01-01"
$ws.Range("J19").ClearFormats()

# Row 20: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.100.0.4 (moved up from row 21)
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.100.0.4"
$ws.Range("B20").Value = 240434
$ws.Range("C20").Value = 421
$ws.Range("D20").Value = 37
$ws.Range("E20").Value = 409
$ws.Range("F20").Value = 240892
$ws.Range("H20").Value = "23.100.0.4"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "2024-11-10"
$ws.Range("J20").ClearFormats()

# Row 21: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.1.3 (moved down from row 20)
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.1.3"
$ws.Range("B21").Value = 151287
$ws.Range("C21").Value = 285
$ws.Range("D21").Value = 75
$ws.Range("E21").Value = 332
$ws.Range("F21").Value = 151647
$ws.Range("H21").Value = "23.80.1.3"
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = "2024-09-03"
$ws.Range("J21").ClearFormats()

# Row 22: Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1 (unchanged position)
$ws.Range("B22").Value = 272039
$ws.Range("C22").Value = 213
$ws.Range("E22").Value = 316
$ws.Range("F22").Value = 272383

# Row 23: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3 (unchanged position)
$ws.Range("B23").Value = 14561
$ws.Range("E23").Value = 59
$ws.Range("F23").Value = 14561

# Row 24: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1 (unchanged position)
$ws.Range("B24").Value = 12018
$ws.Range("E24").Value = 61
$ws.Range("F24").Value = 12018
